$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.093312672893958393
$ws.Range("B1").Value = 0.092953307894433124
$ws.Range("A2").Value = 0.0079272795796541118
$ws.Range("B2").Value = -0.0086854678344554515
$ws.Range("A3").Value = 0.11161918458608255
$ws.Range("B3").Value = -0.11216596182085325
$ws.Range("A4").Value = -0.17183040035058994
$ws.Range("B4").Value = 0.17089547005690164
$ws.Range("A5").Value = -0.1648954702180232
$ws.Range("B5").Value = 0.16301909625413469
$ws.Range("A6").Value = -0.11637028851423414
$ws.Range("B6").Value = 0.11617178088168512
$ws.Range("A7").Value = -0.089212646473312063
$ws.Range("B7").Value = 0.08876013459188492
$ws.Range("A8").Value = -0.068760134796039374
$ws.Range("B8").Value = 0.068369122305917074
$ws.Range("A9").Value = -0.062369122479298156
$ws.Range("B9").Value = 0.06203614621869491
$ws.Range("A10").Value = -0.056036146394937703
$ws.Range("B10").Value = 0.055986254014989356
$ws.Range("A11").Value = -0.051486254187995684
$ws.Range("B11").Value = 0.051404835196628085
$ws.Range("A12").Value = -0.04540483537402551
$ws.Range("B12").Value = 0.045155533660919733
$ws.Range("A13").Value = -0.039155533841342738
$ws.Range("B13").Value = 0.039087503754398512
$ws.Range("A14").Value = -0.027087503950284919
$ws.Range("B14").Value = 0.027054406781561013
$ws.Range("A15").Value = -0.02105440696363825
$ws.Range("B15").Value = 0.0210283331670853
$ws.Range("A16").Value = -0.015028333349822232
$ws.Range("B16").Value = 0.015004764383160785
$ws.Range("A17").Value = -0.0090047645667663545
$ws.Range("B17").Value = 0.008999999808903425
$ws.Range("A18").Value = -0.07843007609409014
$ws.Range("B18").Value = 0.078338824029366094
$ws.Range("A19").Value = -0.06933882419702142
$ws.Range("B19").Value = 0.068653930784832173
$ws.Range("A20").Value = -0.059653930956772072
$ws.Range("B20").Value = 0.059507938055736176
$ws.Range("A21").Value = -0.050507938228614435
$ws.Range("B21").Value = 0.050315594525227159
$ws.Range("A22").Value = -0.093933716918867205
$ws.Range("B22").Value = 0.093624924964810674
$ws.Range("A23").Value = -0.084624925133632622
$ws.Range("B23").Value = 0.084125012960150514
$ws.Range("A24").Value = -0.042125013211694373
$ws.Range("B24").Value = 0.041999999747050154
$ws.Range("A25").Value = -0.040648807906677575
$ws.Range("B25").Value = 0.04060918651581602
$ws.Range("A26").Value = -0.034609186682882154
$ws.Range("B26").Value = 0.034564237965895472
$ws.Range("A27").Value = -0.028564238133275133
$ws.Range("B27").Value = 0.028430778939081414
$ws.Range("A28").Value = -0.022430779107687648
$ws.Range("B28").Value = 0.022352497218418144
$ws.Range("A29").Value = -0.010352497402262628
$ws.Range("B29").Value = 0.01032974602634873
$ws.Range("A30").Value = -0.042163194840339813
$ws.Range("B30").Value = 0.042019725374669026
$ws.Range("A31").Value = -0.02701972556765142
$ws.Range("B31").Value = 0.027000991926239948
$ws.Range("A32").Value = -0.0060009921340329342
$ws.Range("B32").Value = 0.0059999998284183675

$ws.Columns.Item(2).ColumnWidth = 14.7109375

